# Auto-generated script to apply 2023-05-28 crime data updates to column J (2023 totals)
# across the "Citywide Totals", "By Neighborhood", and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 2861
$ws.Cells.Item(3, 10).Value = 2947
$ws.Cells.Item(4, 10).Value = 667
$ws.Cells.Item(5, 10).Value = 232
$ws.Cells.Item(6, 10).Value = 3607
$ws.Cells.Item(7, 10).Value = 10314

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 10).Value = 37
$ws.Cells.Item(7, 10).Value = 120

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(2, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 38

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(6, 10).Value = 96
$ws.Cells.Item(7, 10).Value = 342

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 48
$ws.Cells.Item(6, 10).Value = 39
$ws.Cells.Item(7, 10).Value = 149

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 10).Value = 78
$ws.Cells.Item(3, 10).Value = 152
$ws.Cells.Item(7, 10).Value = 373

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 10).Value = 81
$ws.Cells.Item(3, 10).Value = 76
$ws.Cells.Item(6, 10).Value = 96
$ws.Cells.Item(7, 10).Value = 270

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 10).Value = 81
$ws.Cells.Item(7, 10).Value = 312
$ws.Cells.Item(8, 10).Value = 654
$ws.Cells.Item(13, 10).Value = 14
$ws.Cells.Item(14, 10).Value = 38
$ws.Cells.Item(15, 10).Value = 121
$ws.Cells.Item(19, 10).Value = 326
$ws.Cells.Item(21, 10).Value = 18
$ws.Cells.Item(23, 10).Value = 107
$ws.Cells.Item(24, 10).Value = 31
$ws.Cells.Item(25, 10).Value = 59
$ws.Cells.Item(27, 10).Value = 62
$ws.Cells.Item(29, 10).Value = 586
$ws.Cells.Item(33, 10).Value = 431
$ws.Cells.Item(36, 10).Value = 148
$ws.Cells.Item(37, 10).Value = 342
$ws.Cells.Item(42, 10).Value = 410
$ws.Cells.Item(43, 10).Value = 90
$ws.Cells.Item(47, 10).Value = 81
$ws.Cells.Item(48, 10).Value = 103
$ws.Cells.Item(49, 10).Value = 64
$ws.Cells.Item(51, 10).Value = 139
$ws.Cells.Item(52, 10).Value = 270
$ws.Cells.Item(53, 10).Value = 102
$ws.Cells.Item(54, 10).Value = 201
$ws.Cells.Item(55, 10).Value = 130
$ws.Cells.Item(65, 10).Value = 270
$ws.Cells.Item(67, 10).Value = 373
$ws.Cells.Item(69, 10).Value = 24
$ws.Cells.Item(72, 10).Value = 37
$ws.Cells.Item(73, 10).Value = 95
$ws.Cells.Item(75, 10).Value = 31
$ws.Cells.Item(76, 10).Value = 148
$ws.Cells.Item(77, 10).Value = 88
$ws.Cells.Item(78, 10).Value = 135
$ws.Cells.Item(79, 10).Value = 305
$ws.Cells.Item(82, 10).Value = 12
$ws.Cells.Item(83, 10).Value = 242
$ws.Cells.Item(85, 10).Value = 474
$ws.Cells.Item(90, 10).Value = 115
$ws.Cells.Item(91, 10).Value = 119
$ws.Cells.Item(92, 10).Value = 30
$ws.Cells.Item(93, 10).Value = 51
$ws.Cells.Item(94, 10).Value = 93
$ws.Cells.Item(95, 10).Value = 167
$ws.Cells.Item(96, 10).Value = 120
$ws.Cells.Item(99, 10).Value = 149
$ws.Cells.Item(101, 10).Value = 10314

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 10).Value = 71
$ws.Cells.Item(3, 10).Value = 86
$ws.Cells.Item(7, 10).Value = 242

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 10).Value = 49
$ws.Cells.Item(6, 10).Value = 44
$ws.Cells.Item(7, 10).Value = 167

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 111
$ws.Cells.Item(3, 10).Value = 133
$ws.Cells.Item(5, 10).Value = 17
$ws.Cells.Item(6, 10).Value = 148
$ws.Cells.Item(7, 10).Value = 431

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(6, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 64

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 10).Value = 52
$ws.Cells.Item(7, 10).Value = 201

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 171
$ws.Cells.Item(4, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 586

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 10).Value = 91
$ws.Cells.Item(4, 10).Value = 16
$ws.Cells.Item(6, 10).Value = 126
$ws.Cells.Item(7, 10).Value = 326

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(6, 10).Value = 50
$ws.Cells.Item(7, 10).Value = 103

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 10).Value = 31
$ws.Cells.Item(6, 10).Value = 81
$ws.Cells.Item(7, 10).Value = 148

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 112
$ws.Cells.Item(3, 10).Value = 181
$ws.Cells.Item(6, 10).Value = 135
$ws.Cells.Item(7, 10).Value = 474

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 10).Value = 81
$ws.Cells.Item(3, 10).Value = 87
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 410

$ws = $wb.Worksheets.Item('Boystown')
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(6, 10).Value = 14

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 10).Value = 34
$ws.Cells.Item(7, 10).Value = 135

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 10).Value = 33
$ws.Cells.Item(3, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 130

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(3, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 31

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(6, 10).Value = 23
$ws.Cells.Item(7, 10).Value = 107

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Cells.Item(2, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 24

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(7, 10).Value = 119

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(6, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 18

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 10).Value = 86
$ws.Cells.Item(3, 10).Value = 112
$ws.Cells.Item(7, 10).Value = 305

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 10).Value = 55
$ws.Cells.Item(7, 10).Value = 148

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 10).Value = 11
$ws.Cells.Item(6, 10).Value = 16
$ws.Cells.Item(7, 10).Value = 51

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 10).Value = 61
$ws.Cells.Item(6, 10).Value = 121
$ws.Cells.Item(7, 10).Value = 270

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 93

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 59

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 81

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 121

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 95

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(3, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 81

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(3, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 30

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(3, 10).Value = 208
$ws.Cells.Item(6, 10).Value = 191
$ws.Cells.Item(7, 10).Value = 654

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(3, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 62

$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(3, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 31

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 10).Value = 38
$ws.Cells.Item(3, 10).Value = 34
$ws.Cells.Item(5, 10).Value = 4
$ws.Cells.Item(7, 10).Value = 115

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 10).Value = 36
$ws.Cells.Item(6, 10).Value = 43
$ws.Cells.Item(7, 10).Value = 139

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(4, 10).Value = 9
$ws.Cells.Item(7, 10).Value = 90

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 10).Value = 60
$ws.Cells.Item(7, 10).Value = 102

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(2, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 37

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Cells.Item(5, 10).Value = 6
$ws.Cells.Item(6, 10).Value = 12

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(2, 10).Value = 28
$ws.Cells.Item(7, 10).Value = 88

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 102
$ws.Cells.Item(3, 10).Value = 89
$ws.Cells.Item(6, 10).Value = 107
$ws.Cells.Item(7, 10).Value = 312

Write-Host "Applied 2023-05-28 updates to $($wb.Worksheets.Count) sheet references across 49 worksheets."